# update read file excel in HappyTestCase.java
# Populate the "Element" sheet with Selenium-style locators used by the
# HappyTestCase test, and make it the active sheet/tab (it was previously
# the "data" sheet that was active).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("data")
$ws2 = $wb.Worksheets.Item("Element")

# --- Fill in the "Element" sheet -----------------------------------------
# Row 1 = friendly field names, Row 2 = the matching Selenium locators.
$ws2.Range("A1").Value = "error"
$ws2.Range("A2").Value = "xpath=./ancestor::div[contains(@class,'mktoFieldWrap')]/div[contains(@class,'mktoError')]"
$ws2.Range("B1").Value = "checkbox"
$ws2.Range("B2").Value = "xpath=//div[contains(@class,'mktoCheckboxList')]//input"
$ws2.Range("C1").Value = "email"
$ws2.Range("C2").Value = "id=Email"
$ws2.Range("D1").Value = "btnSubmit"
$ws2.Range("D2").Value = "css=button.mktoButton"
$ws2.Range("E1").Value = "FirstName"
$ws2.Range("F1").Value = "LastName"
$ws2.Range("E2").Value = "id=FirstName"
$ws2.Range("F2").Value = "id=LastName"
$ws2.Range("G1").Value = "Phone"
$ws2.Range("G2").Value = "id=Phone"
$ws2.Range("H1").Value = "Country"
$ws2.Range("H2").Value = "id=Country"
$ws2.Range("I1").Value = "Company"
$ws2.Range("I2").Value = "id=Company"
$ws2.Range("J2").Value = "id=Solution_Interest__c"
$ws2.Range("J1").Value = "selectInterest"
$ws2.Range("K2").Value = "id=Sales_Contact_Comments__c"
$ws2.Range("K1").Value = "areaComment"

# Formatting to match the real workbook: wrapped text, a narrower first
# column, and row heights sized for the wrapped locator strings.
$ws2.Range("A1:K2").WrapText = $true
$ws2.Columns.Item(1).ColumnWidth = 15.14
$ws2.Rows.Item(1).RowHeight = 34
$ws2.Rows.Item(2).RowHeight = 85
$ws2.PageSetup.Orientation = 1

# --- Switch the active tab from "data" to "Element" -----------------------
# Selecting a cell on "data" first records its old cursor position, then
# selecting a cell on "Element" makes it the active/selected sheet/tab.
$ws1.Range("J6").Select()
$ws2.Range("G28").Select()
